$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared string "user_mobile" as header for a new column V (22) ---
$ws.Cells.Item(1, 22).Value = "user_mobile"

# --- New column V (22) data: user_mobile values for rows 2-17 ---
$ws.Cells.Item(2, 22).Value  = 9874563218
$ws.Cells.Item(3, 22).Value  = 9874563218
$ws.Cells.Item(4, 22).Value  = 9874563218
$ws.Cells.Item(5, 22).Value  = 9874563218
$ws.Cells.Item(6, 22).Value  = 9874563218
$ws.Cells.Item(7, 22).Value  = 9874563218
$ws.Cells.Item(8, 22).Value  = 9874563218
$ws.Cells.Item(9, 22).Value  = 9874563218
$ws.Cells.Item(10, 22).Value = 5555555550
$ws.Cells.Item(11, 22).Value = 9874563218
$ws.Cells.Item(12, 22).Value = 9874563218
$ws.Cells.Item(13, 22).Value = 9874563218
$ws.Cells.Item(14, 22).Value = 9874563218
$ws.Cells.Item(15, 22).Value = 9874563218
$ws.Cells.Item(16, 22).Value = 9874653219
$ws.Cells.Item(17, 22).Value = 9874563218

# --- Formatting: copy cell styles from equivalent existing cells so the new
#     column matches the look of the rest of the sheet (style "3" for the
#     regular rows, style "6" for the two rows that already use the
#     highlighted/alternate style on their other text columns) ---
$ws.Range("C1").Copy()
$ws.Range("V1").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("V2:V9").PasteSpecial(-4122)

$ws.Range("E10").Copy()
$ws.Range("V10").PasteSpecial(-4122)

$ws.Range("C11").Copy()
$ws.Range("V11:V15").PasteSpecial(-4122)

$ws.Range("B16").Copy()
$ws.Range("V16").PasteSpecial(-4122)

$ws.Range("C17").Copy()
$ws.Range("V17").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Column width for the new column (bestFit-style width matching column E) ---
$ws.Columns.Item(22).ColumnWidth = 21.83

# --- Selection moved from E23 to F17 ---
$ws.Range("F17").Select()
